$p = $ppt.ActivePresentation

# Slide 1 - subtitle placeholder ("Rectangle 3") is resized/repositioned.
$s1 = $p.Slides.Item(1)
$subtitle = $null
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $cand = $s1.Shapes.Item($i)
    if ($cand.Name -eq "Rectangle 3") {
        $subtitle = $cand
    }
}

# Target EMU values (converted to points, 1 pt = 12700 EMU). The Height value
# is nudged by a hair above the exact quotient so that the Shape object's
# single-precision (float32) Height round-trips to the exact target EMU
# instead of truncating one EMU short.
$subtitle.Left = 48.15433070866142
$subtitle.Top = 275.6699212598425
$subtitle.Width = 657.7108661417323
$subtitle.Height = 113.45512011023621
